$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) retains text formatting for numeric-looking values
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "51.101.49"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").Value = "2.961.33"
$ws.Range("E3").Value = "  +0.89%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "380.13"
$ws.Range("E5").Value = "  +1.43%  "
$ws.Range("D6").Value = "102.28"
$ws.Range("E6").Value = "  +0.65%  "
$ws.Range("D7").Value = "0.545"
$ws.Range("E7").Value = "  +1.89%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "0.590"
$ws.Range("E9").Value = "  +1.56%  "
$ws.Range("D10").Value = "36.67"
$ws.Range("E10").Value = "  +1.14%  "
$ws.Range("E11").Value = "  -0.97%  "
$ws.Range("D12").Value = "0.0854"
$ws.Range("E12").Value = "  +2.22%  "
$ws.Range("D13").Value = "3.420.01"
$ws.Range("E13").Value = "  +0.86%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "18.32"
$ws.Range("E14").Value = "  +2.22%  "
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "7.77"
$ws.Range("E15").Value = "  +6.07%  "
$ws.Range("D16").Value = "11.92"
$ws.Range("E16").Value = "  +66.68%  "
$ws.Range("D17").Value = "2.952.34"
$ws.Range("E17").Value = "  +1.54%  "
$ws.Range("D18").Value = "1.00"
$ws.Range("E18").Value = "  +2.64%  "
$ws.Range("D19").Value = "51.148.07"
$ws.Range("E19").Value = "  +0.43%  "
$ws.Range("D20").Value = "3.10"
$ws.Range("E20").Value = "  -1.29%  "
$ws.Range("D21").Value = "12.42"
$ws.Range("E21").Value = "  -0.63%  "
$ws.Range("D22").Value = "0.0₃0962"
$ws.Range("E22").Value = "  +0.78%  "
$ws.Range("D23").Value = "70.04"
$ws.Range("E23").Value = "  +2.59%  "
$ws.Range("D24").Value = "3.29"
$ws.Range("E24").Value = "  +14.24%  "
$ws.Range("D25").Value = "267.72"
$ws.Range("E25").Value = "  +1.17%  "
$ws.Range("D26").Value = "7.93"
$ws.Range("E26").Value = "  -1.45%  "
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").Value = "7.16"
$ws.Range("E27").Value = "  -6.97%  "
$ws.Range("B28").Value = "Dai"
$ws.Range("C28").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  +0.06%  "
$ws.Range("E29").Value = "  -0.31%  "
$ws.Range("D30").Value = "25.90"
$ws.Range("E30").Value = "  +1.20%  "
$ws.Range("E31").Value = "  -1.78%  "
$ws.Range("D32").Value = "10.43"
$ws.Range("E32").Value = "  +5.91%  "
$ws.Range("D33").Value = "34.38"
$ws.Range("E33").Value = "  +2.55%  "
$ws.Range("B34").Value = "Toncoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D34").Value = "2.08"
$ws.Range("E34").Value = "  +3.05%  "
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").Value = "51.06"
$ws.Range("E35").Value = "  +0.22%  "
$ws.Range("D36").Value = "0.0436"
$ws.Range("E36").Value = "  -2.89%  "
$ws.Range("E37").Value = "  +0.08%  "
$ws.Range("D38").Value = "3.26"
$ws.Range("E38").Value = "  +9.82%  "
$ws.Range("E39").Value = "  +2.03%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").Value = "1.84"
$ws.Range("E40").Value = "  +3.16%  "
$ws.Range("B41").Value = "Celestia"
$ws.Range("C41").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D41").Value = "16.60"
$ws.Range("E41").Value = "  +1.45%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "2.50"
$ws.Range("E42").Value = "  -1.22%  "
$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").Value = "124.74"
$ws.Range("E43").Value = "  +3.85%  "
$ws.Range("D44").Value = "21.55"
$ws.Range("E44").Value = "  +2.84%  "
$ws.Range("D45").Value = "3.53"
$ws.Range("E45").Value = "  +10.04%  "
$ws.Range("D46").Value = "2.39"
$ws.Range("E46").Value = "  +3.37%  "
$ws.Range("E47").Value = "  -1.22%  "
$ws.Range("D48").Value = "2.052.20"
$ws.Range("E48").Value = "  +4.25%  "
$ws.Range("D49").Value = "0.268"
$ws.Range("E49").Value = "  -5.53%  "
$ws.Range("D50").Value = "0.0320"
$ws.Range("E50").Value = "  -6.70%  "
$ws.Range("D51").Value = "5.41"
$ws.Range("E51").Value = "  +7.68%  "

Write-Output "Applied 113 cell updates"
